# Update automatic: dades i banners [2026-02-15 18:50]
# Refresh DATA_EXTRACCIO timestamps and re-scraped measurement values
# from meteo.cat for the 2026-02-15 daily summary sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dades_Meteo")

# Each entry is the target cell (A1 ref) and its new literal text.
# New values are written as a formula returning the literal string,
# then converted to a static value via copy / paste-special-values so
# Excel does not "smart type" percentages (e.g. "59%") into numbers
# and reformat the cell -- the source cells are plain text.
$updates = @(
    @{Cell="E2"; Text='2026-02-15 18:48:36'},
    @{Cell="E3"; Text='2026-02-15 18:48:38'},
    @{Cell="G3"; Text='194 cm'},
    @{Cell="O3"; Text='-5.4 °C'},
    @{Cell="E4"; Text='2026-02-15 18:48:41'},
    @{Cell="O4"; Text='7.4 °C'},
    @{Cell="E5"; Text='2026-02-15 18:48:43'},
    @{Cell="I5"; Text='3.5 mm'},
    @{Cell="O5"; Text='-4.8 °C'},
    @{Cell="E6"; Text='2026-02-15 18:48:46'},
    @{Cell="H6"; Text='59%'},
    @{Cell="E7"; Text='2026-02-15 18:48:48'},
    @{Cell="E8"; Text='2026-02-15 18:48:51'},
    @{Cell="H8"; Text='63%'},
    @{Cell="O8"; Text='8.0 °C'},
    @{Cell="E9"; Text='2026-02-15 18:48:53'},
    @{Cell="O9"; Text='11.0 °C'},
    @{Cell="E10"; Text='2026-02-15 18:48:56'},
    @{Cell="H10"; Text='69%'},
    @{Cell="O10"; Text='7.5 °C'},
    @{Cell="E11"; Text='2026-02-15 18:48:58'},
    @{Cell="H11"; Text='40%'},
    @{Cell="O11"; Text='7.5 °C'},
    @{Cell="E12"; Text='2026-02-15 18:49:00'},
    @{Cell="H12"; Text='53%'},
    @{Cell="E13"; Text='2026-02-15 18:49:03'},
    @{Cell="H13"; Text='34%'},
    @{Cell="J13"; Text='1015.3 hPa'},
    @{Cell="O13"; Text='6.5 °C'},
    @{Cell="E14"; Text='2026-02-15 18:49:05'},
    @{Cell="E15"; Text='2026-02-15 18:49:08'},
    @{Cell="E16"; Text='2026-02-15 18:49:10'},
    @{Cell="H16"; Text='60%'},
    @{Cell="O16"; Text='-2.1 °C'},
    @{Cell="E17"; Text='2026-02-15 18:49:12'},
    @{Cell="H17"; Text='36%'},
    @{Cell="E18"; Text='2026-02-15 18:49:15'},
    @{Cell="O18"; Text='7.3 °C'},
    @{Cell="E19"; Text='2026-02-15 18:49:18'},
    @{Cell="O19"; Text='3.2 °C'},
    @{Cell="E20"; Text='2026-02-15 18:49:20'},
    @{Cell="H20"; Text='59%'},
    @{Cell="E21"; Text='2026-02-15 18:49:23'},
    @{Cell="H21"; Text='37%'},
    @{Cell="E22"; Text='2026-02-15 18:49:25'},
    @{Cell="E23"; Text='2026-02-15 18:49:28'},
    @{Cell="H23"; Text='64%'},
    @{Cell="I23"; Text='1.7 mm'},
    @{Cell="O23"; Text='-3.8 °C'},
    @{Cell="E24"; Text='2026-02-15 18:49:30'},
    @{Cell="O24"; Text='8.7 °C'},
    @{Cell="E25"; Text='2026-02-15 18:49:33'},
    @{Cell="M25"; Text='2.4 °C 18:27 TU'},
    @{Cell="O25"; Text='-1.8 °C'},
    @{Cell="E26"; Text='2026-02-15 18:49:35'},
    @{Cell="E27"; Text='2026-02-15 18:49:37'},
    @{Cell="E28"; Text='2026-02-15 18:49:40'},
    @{Cell="O28"; Text='6.6 °C'},
    @{Cell="E29"; Text='2026-02-15 18:49:42'},
    @{Cell="E30"; Text='2026-02-15 18:49:45'},
    @{Cell="O30"; Text='9.9 °C'},
    @{Cell="E31"; Text='2026-02-15 18:49:47'},
    @{Cell="J31"; Text='1014.3 hPa'},
    @{Cell="O31"; Text='9.8 °C'},
    @{Cell="E32"; Text='2026-02-15 18:49:50'},
    @{Cell="O32"; Text='3.5 °C'},
    @{Cell="E33"; Text='2026-02-15 18:49:52'},
    @{Cell="J33"; Text='1015.3 hPa'},
    @{Cell="O33"; Text='5.6 °C'},
    @{Cell="E34"; Text='2026-02-15 18:49:54'},
    @{Cell="O34"; Text='0.9 °C'},
    @{Cell="E35"; Text='2026-02-15 18:49:57'},
    @{Cell="O35"; Text='4.0 °C'},
    @{Cell="E36"; Text='2026-02-15 18:49:59'},
    @{Cell="H36"; Text='46%'},
    @{Cell="E37"; Text='2026-02-15 18:50:02'},
    @{Cell="E38"; Text='2026-02-15 18:50:04'},
    @{Cell="H38"; Text='66%'},
    @{Cell="O38"; Text='7.7 °C'},
    @{Cell="E39"; Text='2026-02-15 18:50:07'},
    @{Cell="M39"; Text='1.2 °C 18:05 TU'},
    @{Cell="O39"; Text='-3.1 °C'},
    @{Cell="E40"; Text='2026-02-15 18:50:09'},
    @{Cell="H40"; Text='34%'},
    @{Cell="E41"; Text='2026-02-15 18:50:12'},
    @{Cell="E42"; Text='2026-02-15 18:50:15'},
    @{Cell="H42"; Text='55%'},
    @{Cell="E43"; Text='2026-02-15 18:50:17'},
    @{Cell="E44"; Text='2026-02-15 18:50:19'},
    @{Cell="I44"; Text='1.5 mm'},
    @{Cell="O44"; Text='-4.2 °C'},
    @{Cell="E45"; Text='2026-02-15 18:50:22'},
    @{Cell="I45"; Text='0.4 mm'},
    @{Cell="J45"; Text='1023.5 hPa'},
    @{Cell="E46"; Text='2026-02-15 18:50:25'},
    @{Cell="H46"; Text='51%'},
    @{Cell="K46"; Text='12.4 MJ/m2'}
)

foreach ($u in $updates) {
    $r = $ws.Range($u.Cell)
    $r.Formula = '="' + $u.Text.Replace('"', '""') + '"'
    $r.Copy()
    $r.PasteSpecial(-4163)
}
$excel.CutCopyMode = $false
